# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions refresh). Rows 21 and 22 also swap their
# Avalanche/Dai identities (coin order changed upstream), so those two
# rows update Coin name + Link + Price + Volume together.
#
# Values are assigned with a leading apostrophe so Excel stores them as
# literal text (matching the source data's inline-string cells) instead
# of auto-coercing numeric-looking strings (e.g. "11.80", "0.00001011")
# into floating point numbers, which would silently drop meaningful
# trailing/leading zeros. The Style reset afterwards keeps the cell on
# the workbook's default "Normal" style, since the quote-prefix text
# entry would otherwise tag the cell with an implicit Text number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.226.52"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  -0.42%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'1.828.88"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  -0.70%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = "'  +0.21%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'236.24"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  -1.26%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = "'0.6127"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  -2.46%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = "'1.002"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "'  +0.20%  "
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = "'0.07114"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Value = "'0.2817"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  -2.58%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = "'23.58"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'  -5.58%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'0.07674"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  -0.73%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = "'1.829.96"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'  -0.39%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = "'4.821"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  -2.89%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'0.00001011"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  -1.48%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = "'0.6333"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'  -6.25%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'2.067.23"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  -1.21%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'78.98"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  -3.13%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('E18').Value = "'  -5.82%  "
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = "'29.201.11"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  -0.69%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = "'227.95"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "'  -0.41%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('B21').Value = "'Dai"
$ws.Range('B21').Style = "Normal"
$ws.Range('C21').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C21').Style = "Normal"
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  +0.15%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('B22').Value = "'Avalanche"
$ws.Range('B22').Style = "Normal"
$ws.Range('C22').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('C22').Style = "Normal"
$ws.Range('D22').Value = "'11.80"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "'  -4.10%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = "'7.014"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Value = "'1.002"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'  +0.15%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'155.46"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  -1.65%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('E26').Value = "'  -2.22%  "
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = "'8.047"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  -5.16%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = "'16.62"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'  -4.37%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = "'1.486"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'  +1.66%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = "'0.06367"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  -11.85%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = "'1.453"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  -1.78%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = "'3.828"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  -5.47%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = "'3.807"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  -5.56%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'1.131"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  -0.64%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = "'1.748"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  -4.30%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = "'0.6502"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "'  -6.55%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = "'2.549"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'  -1.12%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = "'2.753"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  -1.97%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = "'1.218.24"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  -1.33%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'6.580"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  -3.33%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = "'  -5.40%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = "'0.9197"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  -0.70%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = "'1.001"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  +0.11%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = "'101.45"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  +0.89%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = "'1.973.39"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  -1.50%  "
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = "'62.99"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  -3.51%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = "'0.00000000117"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'  -2.23%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = "'1.623"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'  -5.33%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = "'8.600"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Value = "'0.4568"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  -0.45%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = "'0.05528"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'  -2.51%  "
$ws.Range('E51').Style = "Normal"
